# verbs.xlsx update
#  - fills in the "Potential Form" column (G) for rows 68-84 that were
#    placeholder "na"
#  - adds a new "Passive Form" column (H)
#  - fixes row 2 (行く) Volitional/Potential values and adds its Passive form
#  - appends three new verb rows (107: 誘う, 108: 答える, 109: 残す)
#
# Cell values are written in the same order the original author entered
# them so that newly-interned shared strings line up with the source
# workbook's ordering; formatting is applied afterwards via copy/paste of
# an existing, correctly-styled cell so fonts match the sheet's existing
# styles exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1) New row 107 - 誘う (to invite)
# ---------------------------------------------------------------------
$ws.Range("A107").Value = "誘う"
$ws.Range("B107").Value = "誘って"
$ws.Range("C107").Value = "誘った"
$ws.Range("D107").Value = "誘わない"
$ws.Range("E107").Value = "誘います"
$ws.Range("F107").Value = "誘おう"
$ws.Range("G107").Value = "誘える"

# ---------------------------------------------------------------------
# 2) Potential Form (column G) for rows 68-84 - replace "na" placeholders
# ---------------------------------------------------------------------
$ws.Range("G68").Value = "渡れる"
$ws.Range("G69").Value = "払える"
$ws.Range("G70").Value = "謝れる"
$ws.Range("G71").Value = "困れる"
$ws.Range("G72").Value = "始まれる"
$ws.Range("G73").Value = "終われる"
$ws.Range("G74").Value = "掛れる"
$ws.Range("G75").Value = "押せる"
$ws.Range("G76").Value = "渡せる"
$ws.Range("G77").Value = "返せる"
$ws.Range("G78").Value = "焼ける"
$ws.Range("G79").Value = "とおれる"
$ws.Range("G80").Value = "かよえる"
$ws.Range("G81").Value = "送れる"
$ws.Range("G82").Value = "上がれる"
$ws.Range("G83").Value = "下ろせる"
$ws.Range("G84").Value = "下がれる"

# ---------------------------------------------------------------------
# 3) New row 108 - 答える (to answer)
# ---------------------------------------------------------------------
$ws.Range("A108").Value = "答える"
$ws.Range("B108").Value = "答えて"
$ws.Range("C108").Value = "答えた"
$ws.Range("D108").Value = "答えない"
$ws.Range("E108").Value = "答えます"
$ws.Range("F108").Value = "答えよう"
$ws.Range("G108").Value = "答えられる"

# ---------------------------------------------------------------------
# 4) New column H - "Passive Form"
# ---------------------------------------------------------------------
$ws.Range("H1").Value = "Passive Form"

# ---------------------------------------------------------------------
# 5) Fix row 2 (行く) - Volitional/Potential were wrong (言おう/言える,
#    leftover from 言う) and add its Passive form
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "行こう"
$ws.Range("G2").Value = "行ける"
$ws.Range("H2").Value = "行かれる"

# ---------------------------------------------------------------------
# 6) New row 109 - 残す (to leave behind)
# ---------------------------------------------------------------------
$ws.Range("A109").Value = "残す"
$ws.Range("B109").Value = "残して"
$ws.Range("C109").Value = "残した"
$ws.Range("D109").Value = "残さない"
$ws.Range("E109").Value = "残します"
$ws.Range("F109").Value = "残そう"
$ws.Range("G109").Value = "残せる"

# =======================================================================
# Formatting - copy formats from existing, correctly-styled cells so the
# new cells pick up the same font/alignment (and thus the same style
# index) as their neighbours.
# =======================================================================

# Rows 107-109 use the same Japanese-font (Yu Gothic) style as row 106.
$ws.Range("A106:G106").Copy() | Out-Null
$ws.Range("A107:G107").PasteSpecial($xlPasteFormats)
$ws.Range("A106:G106").Copy() | Out-Null
$ws.Range("A108:G108").PasteSpecial($xlPasteFormats)
$ws.Range("A106:G106").Copy() | Out-Null
$ws.Range("A109:G109").PasteSpecial($xlPasteFormats)
$ws.Rows.Item(107).RowHeight = 18.75
$ws.Rows.Item(108).RowHeight = 18.75
$ws.Rows.Item(109).RowHeight = 18.75

# G68-G84: switch from the placeholder's "na" font to the Japanese font,
# matching the rest of the Potential Form column. G74 and G77 keep their
# original (placeholder) font, as in the source edit.
$japaneseDonor = $ws.Range("G106")
$japaneseRows = @(68, 69, 70, 71, 72, 73, 75, 76, 78, 79, 80, 81, 82, 83, 84)
foreach ($r in $japaneseRows) {
    $japaneseDonor.Copy() | Out-Null
    $ws.Cells.Item($r, 7).PasteSpecial($xlPasteFormats)
}

# New header cell H1 matches the other column headers (G1); H2 matches
# the rest of row 2 (G2).
$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1").PasteSpecial($xlPasteFormats)
$ws.Range("G2").Copy() | Out-Null
$ws.Range("H2").PasteSpecial($xlPasteFormats)

# Give column H a sensible width, matching the other text columns.
$ws.Columns.Item(8).ColumnWidth = 43.5

# Clear clipboard/marching-ants selection artifact and leave the
# selection where the author's edit session ended up.
$excel.CutCopyMode = $false
$ws.Range("G109").Select()
